# Insert two new data rows just above the current row 388. This pushes the
# existing rows 388-411 down to 390-413 (matching the target dimension
# A1:R413) and leaves two blank rows at 388-389 (inheriting the date-style
# formatting from the surrounding rows) ready to be populated with the new
# records.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("388:389").Insert()

# New row 388
$ws.Cells.Item(388, 1).Value = 5
$ws.Cells.Item(388, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(388, 3).Value = "Maule"
$ws.Cells.Item(388, 4).Value = 44516
$ws.Cells.Item(388, 5).Value = 7
$ws.Cells.Item(388, 6).Value = 100112020
$ws.Cells.Item(388, 7).Value = "Tomate"
$ws.Cells.Item(388, 8).Value = "Larga vida"
$ws.Cells.Item(388, 9).Value = "Primera"
$ws.Cells.Item(388, 10).Value = 3500
$ws.Cells.Item(388, 11).Value = 4000
$ws.Cells.Item(388, 12).Value = 4000
$ws.Cells.Item(388, 13).Value = 4000
$ws.Cells.Item(388, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(388, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(388, 16).Value = 400
$ws.Cells.Item(388, 17).Value = 10
$ws.Cells.Item(388, 18).Value = "Hortaliza"

# New row 389
$ws.Cells.Item(389, 1).Value = 5
$ws.Cells.Item(389, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(389, 3).Value = "Maule"
$ws.Cells.Item(389, 4).Value = 44516
$ws.Cells.Item(389, 5).Value = 7
$ws.Cells.Item(389, 6).Value = 100112020
$ws.Cells.Item(389, 7).Value = "Tomate"
$ws.Cells.Item(389, 8).Value = "Larga vida"
$ws.Cells.Item(389, 9).Value = "Primera"
$ws.Cells.Item(389, 10).Value = 3500
$ws.Cells.Item(389, 11).Value = 7000
$ws.Cells.Item(389, 12).Value = 7000
$ws.Cells.Item(389, 13).Value = 7000
$ws.Cells.Item(389, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(389, 15).Value = "Región del Maule"
$ws.Cells.Item(389, 16).Value = 467
$ws.Cells.Item(389, 17).Value = 15
$ws.Cells.Item(389, 18).Value = "Hortaliza"
